$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "long edge"
$ws.Range("B1").Value = "length"
$ws.Range("C1").Value = "width"
$ws.Range("D1").Value = "short edge"

$ws.Range("C2").Select()
